# Replace ':' with '=' as the column/value separator in the "antecedents"
# (column A) and "consequents" (column B) text for every data row (2-34).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = 34
}

for ($r = 2; $r -le $lastRow; $r++) {
    foreach ($c in 1, 2) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value()
        if ($val -ne $null -and $val.GetType().Name -eq "String" -and $val.Contains(":")) {
            $cell.Value = $val.Replace(":", "=")
        }
    }
}
